# implemented Trace Viewer for Testing Artifacts
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test case identifiers / statuses to the new "trace viewer" style
# names, and collapse the result set from 6 test rows down to 4.

$ws.Range("A2").Value = "TC-TESTS.LOGINTEST2.TESTVALIDLOGIN3 [3.0, SECRET_SAUCE, TESTCASE]"
$ws.Range("B2").Value = "testValidLogin3"
$ws.Range("C2").Value = "PASS"

$ws.Range("A3").Value = "TC-TESTS.LOGINTEST.TESTVALIDLOGIN [3.0, SECRET_SAUCE, TESTCASE]"
$ws.Range("B3").Value = "testValidLogin"
$ws.Range("C3").Value = "PASS"

$ws.Range("A4").Value = "TC-TESTS.LOGINTEST2.TESTVALIDLOGIN3 [4.0, SECRET_SAUCE, TESTCASE]"
$ws.Range("B4").Value = "testValidLogin3"
$ws.Range("C4").Value = "PASS"

$ws.Range("A5").Value = "TC-TESTS.LOGINTEST2.TESTVALIDLOGIN3 [5.0, SECRET_SAUCE, TESTCASE]"
$ws.Range("B5").Value = "testValidLogin3"
$ws.Range("C5").Value = "PASS"

# Remove the now-obsolete last two rows of results
$ws.Range("A7:C7").EntireRow.Delete()
$ws.Range("A6:C6").EntireRow.Delete()
